# Generate Report for Handoff
# Renames the source markdown file from f046273b-eaf2-4f5e-bb45-fbe7658c530d
# to d1cff685-bc69-47cd-bc2a-a1a10b9bfebb, regenerates the handoff xliff
# filenames/timestamps, and clears the "already handed back" info (the
# new handoff has not come back from localization yet).

$wb = $excel.ActiveWorkbook

$oldId = "f046273b-eaf2-4f5e-bb45-fbe7658c530d"
$newId = "d1cff685-bc69-47cd-bc2a-a1a10b9bfebb"
$oldHash = "8807652c59e58ec4b71b671a5e306c5f743b2ab0"
$newHash = "18b7408ff9ffbd43904622302e6be40c5660e4c3"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "$newId.md"
$ov.Range("B2").Value = "e2e\$newId.md"
foreach ($h in $ov.Hyperlinks) {
    if ($h.Range.Address() -eq "`$B`$2") {
        $h.TextToDisplay = "e2e\$newId.md"
    }
}
$ov.Range("G2").Value = "2016-08-12 05:02:45"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newId.md"
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = "$newId.md"
    }
}

$zh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-12 05:02:39"

# Remove the "Latest Target File" hyperlink/value - no handback yet.
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq "`$I`$2") {
        $h.Delete()
    }
}
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"

# Latest Handback File / DateTime - cleared (not handed back yet).
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

# Column widths re-autofit now that I/J no longer hold long filenames.
$zh.Columns.Item(9).AutoFit()
$zh.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newId.md"
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.TextToDisplay = "$newId.md"
    }
}

$de.Range("G2").Value = "$newId.$newHash.de-de.xlf"

# Remove the "Latest Target File" hyperlink/value - no handback yet.
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq "`$I`$2") {
        $h.Delete()
    }
}
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"

# Latest Handback File / DateTime - cleared (not handed back yet).
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

# Column widths re-autofit now that I/J no longer hold long filenames.
$de.Columns.Item(9).AutoFit()
$de.Columns.Item(10).AutoFit()
